$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report date (B1) ---
$ws.Range("B1").Value = 45658

# --- Finished goods rows (existing rows 4-21, values updated) ---
$ws.Range("B4").Value = 53144.26
$ws.Range("C4").Value = 257.17490468396772
$ws.Range("D4").Value = 13667370

$ws.Range("B5").Value = 62641.39
$ws.Range("C5").Value = 258.28611881058197
$ws.Range("D5").Value = 16179401.5

$ws.Range("B7").Value = 120.5
$ws.Range("C7").Value = 348.54771784232366
$ws.Range("D7").Value = 42000

# Row 8 (MONOFILAMENT FABRIC  HAPPA): B/D cleared, C stays 0
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 0
$ws.Range("D8").ClearContents()

$ws.Range("B9").Value = 12275
$ws.Range("C9").Value = 132.37270875763747
$ws.Range("D9").Value = 1624875

$ws.Range("B10").Value = 128181.15
$ws.Range("C10").Value = 245.8524244789503
$ws.Range("D10").Value = 31513646.5

$ws.Range("B11").Value = 1203.25
$ws.Range("C11").Value = 275.36796177020568
$ws.Range("D11").Value = 331336.5

$ws.Range("B12").Value = 806
$ws.Range("C12").Value = 219.96153846153845
$ws.Range("D12").Value = 177289

$ws.Range("B13").Value = 7455.4
$ws.Range("C13").Value = 165.40229632212893
$ws.Range("D13").Value = 1233140.28

# Row 15 (Knitted Fabric): B/D cleared, C stays 0
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 0
$ws.Range("D15").ClearContents()

$ws.Range("B16").Value = 300
$ws.Range("C16").Value = 180
$ws.Range("D16").Value = 54000

$ws.Range("B17").Value = 9764.65
$ws.Range("C17").Value = 183.9047769249282
$ws.Range("D17").Value = 1795765.78

$ws.Range("B20").Value = 5060
$ws.Range("D20").Value = 50600

# --- Insert a new row at 22 for "Master Batch" (shifts old rows 22-39 down to 23-40) ---
$ws.Rows(22).Insert()

$ws.Range("A22").Value = "Master Batch"
$ws.Range("B22").Value = 205
$ws.Range("C22").Value = 268.6829268292683
$ws.Range("D22").Value = 55080

# --- Row 23: Raw Material (was row 22 pre-insert) ---
$ws.Range("B23").Value = 39000
$ws.Range("C23").Value = 95.65384615384616
$ws.Range("D23").Value = 3730500

# --- Row 24 (no label) ---
$ws.Range("B24").Value = 44265
$ws.Range("C24").Value = 86.66395572122444
$ws.Range("D24").Value = 3836180

# --- Row 25: Grand total: (now carries B/C/D, previously empty of those) ---
$ws.Range("B25").Value = 182210.8
$ws.Range("C25").Value = 203.86054108757551
$ws.Range("D25").Value = 37145592.280000001

# --- Row 26: Other Income -- B & C cleared, only D remains ---
$ws.Range("B26").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("D26").Value = 28147.32

# --- Row 27 (no label) -- B cleared, only D remains ---
$ws.Range("B27").ClearContents()
$ws.Range("D27").Value = 37094992

# --- Row 28: Gross sales ---
$ws.Range("D28").Value = 39426761.7

# --- Row 29: tax (new D value) ---
$ws.Range("D29").Value = 2514816.69

# --- Row 30: TCS (new D value) ---
$ws.Range("D30").Value = 13012.73

# --- Row 31 (no label) ---
$ws.Range("D31").Value = 2527829.42

# --- Row 32 (no label, previously had A "Discount") ---
$ws.Range("A32").ClearContents()
$ws.Range("D32").Value = 36898932.280000001

# --- Row 33: Discount (new position) ---
$ws.Range("A33").Value = "Discount"
$ws.Range("D33").Value = 246660

# --- Row 34 (no label, previously had A "Credit Note") ---
$ws.Range("A34").ClearContents()
$ws.Range("D34").Value = 37145592.280000001

# --- Row 35: Credit Note (new position) ---
$ws.Range("A35").Value = "Credit Note"
$ws.Range("D35").Value = 0

# --- Row 36 (no label, new empty D value row) ---
$ws.Range("D36").Value = 0

# --- Row 37: PAL I- FINAL SALES +Less Waste & Discount (shifted from 36) ---
$ws.Range("A37").Value = "PAL I- FINAL SALES +Less Waste  & Discount"
$ws.Range("D37").Value = 36848332

# --- Row 40: RM Purchase for sales (shifted from 39) ---
$ws.Range("A40").Value = "RM Purchase for sales"
$ws.Range("B40").Value = 39000
$ws.Range("C40").Value = 96.36
$ws.Range("D40").Value = 3758040
